# Apply cryptos list update (values refreshed by the GitHub Actions scraper).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D column holds plain numeric-looking text (prices). Force text format first so
# Excel does not reinterpret values like "20.80" or "0.9000" as numbers and drop
# the significant trailing zeros / formatting that the scraped text relies on.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D18", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "27.815.87"
$ws.Range("E2").Value = "  -0.95%  "

# Row 3
$ws.Range("D3").Value = "1.903.00"
$ws.Range("E3").Value = "  -0.34%  "

# Row 4
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "313.15"
$ws.Range("E5").Value = "  -0.87%  "

# Row 6
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -0.03%  "

# Row 7
$ws.Range("D7").Value = "0.5022"
$ws.Range("E7").Value = "  +3.96%  "

# Row 8
$ws.Range("D8").Value = "0.3811"
$ws.Range("E8").Value = "  -0.31%  "

# Row 9
$ws.Range("D9").Value = "0.07286"
$ws.Range("E9").Value = "  -0.94%  "

# Row 10
$ws.Range("D10").Value = "0.9086"
$ws.Range("E10").Value = "  -2.77%  "

# Row 11
$ws.Range("D11").Value = "20.80"
$ws.Range("E11").Value = "  -0.04%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.929.08"
$ws.Range("E12").Value = "  +1.00%  "

# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "0.07653"
$ws.Range("E13").Value = "  -1.91%  "

# Row 14
$ws.Range("D14").Value = "5.479"
$ws.Range("E14").Value = "  -0.54%  "

# Row 15
$ws.Range("D15").Value = "6.587"
$ws.Range("E15").Value = "  -0.58%  "

# Row 16
$ws.Range("D16").Value = "91.29"
$ws.Range("E16").Value = "  -0.17%  "

# Row 17
$ws.Range("E17").Value = "  -0.18%  "

# Row 18
$ws.Range("D18").Value = "0.000008691"
$ws.Range("E18").Value = "  -1.53%  "

# Row 19
$ws.Range("E19").Value = "  -0.19%  "

# Row 20
$ws.Range("D20").Value = "27.835.21"
$ws.Range("E20").Value = "  -0.94%  "

# Row 21
$ws.Range("D21").Value = "14.51"
$ws.Range("E21").Value = "  -2.31%  "

# Row 22
$ws.Range("D22").Value = "5.158"
$ws.Range("E22").Value = "  +0.22%  "

# Row 23
$ws.Range("D23").Value = "10.82"
$ws.Range("E23").Value = "  -0.66%  "

# Row 24
$ws.Range("D24").Value = "154.15"
$ws.Range("E24").Value = "  -1.80%  "

# Row 25
$ws.Range("D25").Value = "1.849"
$ws.Range("E25").Value = "  -3.95%  "

# Row 26
$ws.Range("D26").Value = "2.233"
$ws.Range("E26").Value = "  +6.09%  "

# Row 27
$ws.Range("D27").Value = "18.37"
$ws.Range("E27").Value = "  -1.10%  "

# Row 28
$ws.Range("D28").Value = "115.27"
$ws.Range("E28").Value = "  -0.91%  "

# Row 29
$ws.Range("D29").Value = "4.910"
$ws.Range("E29").Value = "  -0.90%  "

# Row 30
$ws.Range("D30").Value = "0.08972"
$ws.Range("E30").Value = "  +0.65%  "

# Row 31
$ws.Range("D31").Value = "3.207"
$ws.Range("E31").Value = "  -3.90%  "

# Row 32
$ws.Range("D32").Value = "1.231"
$ws.Range("E32").Value = "  -1.57%  "

# Row 33
$ws.Range("D33").Value = "0.7657"
$ws.Range("E33").Value = "  -0.08%  "

# Row 34
$ws.Range("D34").Value = "4.633"
$ws.Range("E34").Value = "  -1.13%  "

# Row 35
$ws.Range("D35").Value = "0.02058"
$ws.Range("E35").Value = "  +0.43%  "

# Row 36
$ws.Range("D36").Value = "2.549"
$ws.Range("E36").Value = "  -2.33%  "

# Row 37
$ws.Range("D37").Value = "1.099"
$ws.Range("E37").Value = "  -0.25%  "

# Row 38
$ws.Range("D38").Value = "0.5533"
$ws.Range("E38").Value = "  +0.84%  "

# Row 39
$ws.Range("D39").Value = "3.016"
$ws.Range("E39").Value = "  +1.21%  "

# Row 40
$ws.Range("D40").Value = "0.05254"
$ws.Range("E40").Value = "  -0.82%  "

# Row 41
$ws.Range("D41").Value = "6.971"
$ws.Range("E41").Value = "  -0.74%  "

# Row 42
$ws.Range("D42").Value = "8.496"
$ws.Range("E42").Value = "  +0.50%  "

# Row 43
$ws.Range("D43").Value = "0.1523"
$ws.Range("E43").Value = "  +0.01%  "

# Row 44
$ws.Range("D44").Value = "110.95"
$ws.Range("E44").Value = "  +3.50%  "

# Row 45
$ws.Range("D45").Value = "10.59"
$ws.Range("E45").Value = "  -0.77%  "

# Row 46
$ws.Range("D46").Value = "0.4783"
$ws.Range("E46").Value = "  -1.00%  "

# Row 47
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").Value = "  -0.02%  "

# Row 48
$ws.Range("D48").Value = "1.631"
$ws.Range("E48").Value = "  -1.51%  "

# Row 49
$ws.Range("D49").Value = "67.14"
$ws.Range("E49").Value = "  -1.78%  "

# Row 50
$ws.Range("D50").Value = "0.06074"
$ws.Range("E50").Value = "  -0.54%  "

# Row 51
$ws.Range("D51").Value = "0.9000"
$ws.Range("E51").Value = "  -0.25%  "
